$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2474.75
$ws.Range("J43").Value = 2324.75
$ws.Range("L43").Value = 2324.75
$ws.Range("N43").Value = -2462.75
$ws.Range("H98").Value = 1856.2858
$ws.Range("I98").Value = 1750.6666
$ws.Range("J98").Value = 2490
$ws.Range("K98").Value = 1750.6666
$ws.Range("L98").Value = 2490
$ws.Range("M98").Value = -252.6666
$ws.Range("N98").Value = -5486
$ws.Range("H111").Value = 795.875
$ws.Range("I111").Value = 741.2
$ws.Range("K111").Value = 2223.6
$ws.Range("M111").Value = 843.3999999999996
$ws.Range("H122").Value = 1856.2858
$ws.Range("I122").Value = 1750.6666
$ws.Range("J122").Value = 2490
$ws.Range("K122").Value = 5251.9998
$ws.Range("L122").Value = 7470
$ws.Range("M122").Value = -2801.9998
$ws.Range("N122").Value = -12370
$ws.Range("H127").Value = 1655.4286
$ws.Range("I127").Value = 1431.3334
$ws.Range("J127").Value = 3000
$ws.Range("K127").Value = 4294.0002
$ws.Range("L127").Value = 9000
$ws.Range("M127").Value = 665.9997999999996
$ws.Range("N127").Value = -18920
$ws.Range("H132").Value = 3081.9092
$ws.Range("I132").Value = 3310.1
$ws.Range("K132").Value = 9930.299999999999
$ws.Range("M132").Value = -7400.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 820
$ws.Range("I97").Value = 800.9
$ws.Range("K97").Value = 800.9
$ws.Range("M97").Value = -304.9
$ws.Range("H110").Value = 3649
$ws.Range("I110").Value = 1850
$ws.Range("K110").Value = 1850
$ws.Range("M110").Value = 195

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4689.8887
$ws.Range("I99").Value = 5365.7144
$ws.Range("K99").Value = 5365.7144
$ws.Range("M99").Value = -3867.7144
$ws.Range("H107").Value = 1097.4615
$ws.Range("I107").Value = 976.7
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 976.7
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 943.3
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 3358
$ws.Range("I33").Value = 1697.5
$ws.Range("K33").Value = 1697.5
$ws.Range("M33").Value = -1318.5
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H58").Value = 2745.2632
$ws.Range("I58").Value = 1399.5
$ws.Range("J58").Value = 3724
$ws.Range("K58").Value = 1399.5
$ws.Range("L58").Value = 3724
$ws.Range("M58").Value = -1196.5
$ws.Range("N58").Value = -4130
$ws.Range("H105").Value = 1179.9375
$ws.Range("I105").Value = 930.8182
$ws.Range("K105").Value = 930.8182
$ws.Range("M105").Value = 816.1818
$ws.Range("H122").Value = 749
$ws.Range("I122").Value = 770.8125
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 2312.4375
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 137.5625
$ws.Range("N122").Value = -6100
$ws.Range("H136").Value = 2745.2632
$ws.Range("I136").Value = 1399.5
$ws.Range("J136").Value = 3724
$ws.Range("K136").Value = 4198.5
$ws.Range("L136").Value = 11172
$ws.Range("M136").Value = -1648.5
$ws.Range("N136").Value = -16272

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 939.5
$ws.Range("I16").Value = 939
$ws.Range("K16").Value = 2817
$ws.Range("M16").Value = -2644
$ws.Range("H80").Value = 11333
$ws.Range("J80").Value = 11333
$ws.Range("L80").Value = 33999
$ws.Range("N80").Value = -35871
$ws.Range("H83").Value = 11333
$ws.Range("J83").Value = 11333
$ws.Range("L83").Value = 101997
$ws.Range("N83").Value = -111357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 28958.416
$ws.Range("I26").Value = 19750
$ws.Range("K26").Value = 19750
$ws.Range("M26").Value = -19470
$ws.Range("H50").Value = 28958.416
$ws.Range("I50").Value = 19750
$ws.Range("K50").Value = 19750
$ws.Range("M50").Value = -19252
$ws.Range("H70").Value = 111112450
$ws.Range("I70").Value = 166667660
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 166667660
$ws.Range("L70").Value = 2000
$ws.Range("M70").Value = -166667390
$ws.Range("N70").Value = -2540
$ws.Range("H73").Value = 111112450
$ws.Range("I73").Value = 166667660
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 166667660
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = -166666724
$ws.Range("N73").Value = -3872
$ws.Range("H80").Value = 2961.75
$ws.Range("I80").Value = 2813.4285
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2813.4285
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1815.4285
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 2961.75
$ws.Range("I83").Value = 2813.4285
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 14067.1425
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -9075.1425
$ws.Range("N83").Value = -29984
$ws.Range("H113").Value = 902.6
$ws.Range("I113").Value = 925
$ws.Range("J113").Value = 813
$ws.Range("K113").Value = 925
$ws.Range("L113").Value = 813
$ws.Range("M113").Value = 1245
$ws.Range("N113").Value = -5153

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16445.5
$ws.Range("I7").Value = 13052.75
$ws.Range("J7").Value = 26623.75
$ws.Range("K7").Value = 13052.75
$ws.Range("L7").Value = 26623.75
$ws.Range("M7").Value = -12940.75
$ws.Range("N7").Value = -26847.75
$ws.Range("H22").Value = 516.3
$ws.Range("I22").Value = 271.375
$ws.Range("K22").Value = 271.375
$ws.Range("M22").Value = 23.625
$ws.Range("H27").Value = 516.3
$ws.Range("I27").Value = 271.375
$ws.Range("K27").Value = 271.375
$ws.Range("M27").Value = -164.375
$ws.Range("H40").Value = 6752.5
$ws.Range("I40").Value = 6752.5
$ws.Range("K40").Value = 6752.5
$ws.Range("M40").Value = -6616.5
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H55").Value = 239.1
$ws.Range("J55").Value = 267.8
$ws.Range("L55").Value = 267.8
$ws.Range("N55").Value = -613.8
$ws.Range("H126").Value = 16445.5
$ws.Range("I126").Value = 13052.75
$ws.Range("J126").Value = 26623.75
$ws.Range("K126").Value = 39158.25
$ws.Range("L126").Value = 79871.25
$ws.Range("M126").Value = -36688.25
$ws.Range("N126").Value = -84811.25
$ws.Range("H136").Value = 6287.6665
$ws.Range("I136").Value = 5519.2
$ws.Range("J136").Value = 7248.25
$ws.Range("K136").Value = 16557.6
$ws.Range("L136").Value = 21744.75
$ws.Range("M136").Value = -14007.6
$ws.Range("N136").Value = -26844.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360
$ws.Range("H107").Value = 2126.3044
$ws.Range("I107").Value = 2512.8462
$ws.Range("K107").Value = 7538.5386
$ws.Range("M107").Value = -5618.5386
$ws.Range("H126").Value = 78834.336
$ws.Range("I126").Value = 67002
$ws.Range("J126").Value = 84750.5
$ws.Range("K126").Value = 201006
$ws.Range("L126").Value = 254251.5
$ws.Range("M126").Value = -198536
$ws.Range("N126").Value = -259191.5
$ws.Range("H132").Value = 2319.923
$ws.Range("I132").Value = 2319.923
$ws.Range("K132").Value = 6959.768999999999
$ws.Range("M132").Value = -4429.768999999999
